$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.00581"
$ws.Range("H2").Value = [double]"3.01743"
$ws.Range("I2").Value = [double]"0.003799625168827527"
$ws.Range("J2").Value = [double]"0.003799625168827527"
$ws.Range("M2").Value = [double]"0.02270466666666667"
$ws.Range("N2").Value = [double]"0.06811400000000001"
$ws.Range("O2").Value = [double]"0.002206225855740089"
$ws.Range("P2").Value = [double]"0.002206225855740089"
$ws.Range("Q2").Value = [double]"0.02283658078"
$ws.Range("R2").Value = [double]"0.20552922702"
$ws.Range("S2").Value = [double]"8.38283128958809E-06"
$ws.Range("T2").Value = [double]"8.38283128958809E-06"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.00581"
$ws.Range("H3").Value = [double]"3.01743"
$ws.Range("I3").Value = [double]"0.003799625168827527"
$ws.Range("J3").Value = [double]"0.003799625168827527"
$ws.Range("O3").Value = [double]"0.002281111990432972"
$ws.Range("P3").Value = [double]"0.002281111990432972"
$ws.Range("Q3").Value = [double]"0.02361172502"
$ws.Range("R3").Value = [double]"0.21250552518"
$ws.Range("S3").Value = [double]"8.667370531763378E-06"
$ws.Range("T3").Value = [double]"8.667370531763378E-06"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.00581"
$ws.Range("H4").Value = [double]"3.01743"
$ws.Range("I4").Value = [double]"0.003799625168827527"
$ws.Range("J4").Value = [double]"0.003799625168827527"
$ws.Range("M4").Value = [double]"10.24499966666667"
$ws.Range("N4").Value = [double]"30.734999"
$ws.Range("O4").Value = [double]"0.9955126621538269"
$ws.Range("P4").Value = [double]"0.9955126621538269"
$ws.Range("Q4").Value = [double]"10.30452311473"
$ws.Range("R4").Value = [double]"92.74070803257001"
$ws.Range("S4").Value = [double]"0.003782574967006176"
$ws.Range("T4").Value = [double]"0.003782574967006176"
$ws.Range("I5").Value = [double]"0.9594121222074437"
$ws.Range("J5").Value = [double]"0.9594121222074438"
$ws.Range("M5").Value = [double]"0.02270466666666667"
$ws.Range("N5").Value = [double]"0.06811400000000001"
$ws.Range("O5").Value = [double]"0.002206225855740089"
$ws.Range("P5").Value = [double]"0.002206225855740089"
$ws.Range("Q5").Value = [double]"5.766277318576224"
$ws.Range("R5").Value = [double]"51.89649586718601"
$ws.Range("S5").Value = [double]"0.002116679830324532"
$ws.Range("T5").Value = [double]"0.002116679830324532"
$ws.Range("I6").Value = [double]"0.9594121222074437"
$ws.Range("J6").Value = [double]"0.9594121222074438"
$ws.Range("O6").Value = [double]"0.002281111990432972"
$ws.Range("P6").Value = [double]"0.002281111990432972"
$ws.Range("S6").Value = [double]"0.002188526495734144"
$ws.Range("T6").Value = [double]"0.002188526495734144"
$ws.Range("I7").Value = [double]"0.9594121222074437"
$ws.Range("J7").Value = [double]"0.9594121222074438"
$ws.Range("M7").Value = [double]"10.24499966666667"
$ws.Range("N7").Value = [double]"30.734999"
$ws.Range("O7").Value = [double]"0.9955126621538269"
$ws.Range("P7").Value = [double]"0.9955126621538269"
$ws.Range("Q7").Value = [double]"2601.910438678728"
$ws.Range("R7").Value = [double]"23417.19394810855"
$ws.Range("S7").Value = [double]"0.955106915881385"
$ws.Range("T7").Value = [double]"0.9551069158813851"
$ws.Range("G8").Value = [double]"9.336668333333334"
$ws.Range("H8").Value = [double]"28.010005"
$ws.Range("I8").Value = [double]"0.03527091597053946"
$ws.Range("J8").Value = [double]"0.03527091597053946"
$ws.Range("M8").Value = [double]"0.02270466666666667"
$ws.Range("N8").Value = [double]"0.06811400000000001"
$ws.Range("O8").Value = [double]"0.002206225855740089"
$ws.Range("P8").Value = [double]"0.002206225855740089"
$ws.Range("Q8").Value = [double]"0.2119859422855556"
$ws.Range("R8").Value = [double]"1.90787348057"
$ws.Range("S8").Value = [double]"7.781560676984018E-05"
$ws.Range("T8").Value = [double]"7.781560676984018E-05"
$ws.Range("G9").Value = [double]"9.336668333333334"
$ws.Range("H9").Value = [double]"28.010005"
$ws.Range("I9").Value = [double]"0.03527091597053946"
$ws.Range("J9").Value = [double]"0.03527091597053946"
$ws.Range("O9").Value = [double]"0.002281111990432972"
$ws.Range("P9").Value = [double]"0.002281111990432972"
$ws.Range("Q9").Value = [double]"0.2191814013477778"
$ws.Range("R9").Value = [double]"1.97263261213"
$ws.Range("S9").Value = [double]"8.045690933395139E-05"
$ws.Range("T9").Value = [double]"8.045690933395139E-05"
$ws.Range("G10").Value = [double]"9.336668333333334"
$ws.Range("H10").Value = [double]"28.010005"
$ws.Range("I10").Value = [double]"0.03527091597053946"
$ws.Range("J10").Value = [double]"0.03527091597053946"
$ws.Range("M10").Value = [double]"10.24499966666667"
$ws.Range("N10").Value = [double]"30.734999"
$ws.Range("O10").Value = [double]"0.9955126621538269"
$ws.Range("P10").Value = [double]"0.9955126621538269"
$ws.Range("Q10").Value = [double]"95.65416396277723"
$ws.Range("R10").Value = [double]"860.887475664995"
$ws.Range("S10").Value = [double]"0.03511264345443567"
$ws.Range("T10").Value = [double]"0.03511264345443567"
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"0.4016586666666667"
$ws.Range("H11").Value = [double]"1.204976"
$ws.Range("I11").Value = [double]"0.001517336653189343"
$ws.Range("J11").Value = [double]"0.001517336653189343"
$ws.Range("M11").Value = [double]"0.02270466666666667"
$ws.Range("N11").Value = [double]"0.06811400000000001"
$ws.Range("O11").Value = [double]"0.002206225855740089"
$ws.Range("P11").Value = [double]"0.002206225855740089"
$ws.Range("Q11").Value = [double]"0.009119526140444445"
$ws.Range("R11").Value = [double]"0.08207573526400001"
$ws.Range("S11").Value = [double]"3.34758735612846E-06"
$ws.Range("T11").Value = [double]"3.34758735612846E-06"
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.4016586666666667"
$ws.Range("H12").Value = [double]"1.204976"
$ws.Range("I12").Value = [double]"0.001517336653189343"
$ws.Range("J12").Value = [double]"0.001517336653189343"
$ws.Range("O12").Value = [double]"0.002281111990432972"
$ws.Range("P12").Value = [double]"0.002281111990432972"
$ws.Range("Q12").Value = [double]"0.009429071086222223"
$ws.Range("R12").Value = [double]"0.084861639776"
$ws.Range("S12").Value = [double]"3.461214833113646E-06"
$ws.Range("T12").Value = [double]"3.461214833113646E-06"
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.4016586666666667"
$ws.Range("H13").Value = [double]"1.204976"
$ws.Range("I13").Value = [double]"0.001517336653189343"
$ws.Range("J13").Value = [double]"0.001517336653189343"
$ws.Range("M13").Value = [double]"10.24499966666667"
$ws.Range("N13").Value = [double]"30.734999"
$ws.Range("O13").Value = [double]"0.9955126621538269"
$ws.Range("P13").Value = [double]"0.9955126621538269"
$ws.Range("Q13").Value = [double]"4.114992906113778"
$ws.Range("R13").Value = [double]"37.03493615502401"
$ws.Range("S13").Value = [double]"0.0015105278510001"
$ws.Range("T13").Value = [double]"0.0015105278510001"
